# Minpowerinputfile.xlsx edit:
#  - Rename sheet "Band 5" -> "Band5"
#  - Make "Band5" the active/selected tab (was "Band13")
#  - On "Band5", move the selection to G28 (was F2:F4)
# Activating a sheet in real Excel clears tabSelected on the
# previously-active sheet and sets it on the newly activated one, and
# updates the workbook's bookView.activeTab to match - so driving this
# through .Activate()/.Select() reproduces the whole diff in one go.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Band 5")
$ws.Name = "Band5"

$ws.Activate()
$ws.Range("G28").Select()
